# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Only the "Price" (D) and "Volume(1h)" (E) columns move; Coin/Link/rank are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.323.21"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.692.01"
$ws.Range("E3").Value = "  +0.67%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.90"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5411"
$ws.Range("E6").Value = "  +2.90%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.10%  "

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2739"
$ws.Range("E8").Value = "  +1.30%  "

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06456"
$ws.Range("E9").Value = "  -0.44%  "

# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.69"
$ws.Range("E10").Value = "  -1.42%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07666"
$ws.Range("E11").Value = "  +1.68%  "

# Row 12: WrappedEther
$ws.Range("D12").Value = "1.702.98"
$ws.Range("E12").Value = "  +1.54%  "

# Row 13: Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.540"
$ws.Range("E13").Value = "  +0.17%  "

# Row 14: Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5802"

# Row 15: ShibaInu
$ws.Range("E15").Value = "  -1.28%  "

# Row 16: Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.89"
$ws.Range("E16").Value = "  +3.41%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "26.369.51"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18: Uniswap
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.921"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19: Dai
$ws.Range("E19").Value = "  +0.13%  "

# Row 20: Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.88"
$ws.Range("E20").Value = "  -0.02%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.90"
$ws.Range("E21").Value = "  +0.35%  "

# Row 22: Chainlink
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.277"
$ws.Range("E22").Value = "  +1.12%  "

# Row 23: BinanceUSD
$ws.Range("E23").Value = "  +0.11%  "

# Row 24: Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.27"
$ws.Range("E24").Value = "  +2.61%  "

# Row 25: Stellar
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1285"
$ws.Range("E25").Value = "  +3.31%  "

# Row 26: Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.856"
$ws.Range("E26").Value = "  +0.73%  "

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.89"
$ws.Range("E27").Value = "  +0.57%  "

# Row 28: Hedera
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06327"
$ws.Range("E28").Value = "  -2.78%  "

# Row 29: Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.379"
$ws.Range("E29").Value = "  +1.68%  "

# Row 30: PancakeSwap
$ws.Range("E30").Value = "  -0.55%  "

# Row 31: InternetComputer(DFINITY)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.604"
$ws.Range("E31").Value = "  +0.31%  "

# Row 32: Filecoin
$ws.Range("E32").Value = "  -0.34%  "

# Row 33: LidoDAOToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.682"
$ws.Range("E33").Value = "  +1.36%  "

# Row 34: ARBITRUM
$ws.Range("E34").Value = "  -0.04%  "

# Row 35: ImmutableX
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6196"
$ws.Range("E35").Value = "  -0.71%  "

# Row 36: HuobiToken
$ws.Range("E36").Value = "  +0.55%  "

# Row 37: MXToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.772"

# Row 38: VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01658"
$ws.Range("E38").Value = "  +2.02%  "

# Row 39: Maker
$ws.Range("D39").Value = "1.112.90"
$ws.Range("E39").Value = "  +0.08%  "

# Row 40: FraxShare
$ws.Range("E40").Value = "  -5.52%  "

# Row 41: TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8859"
$ws.Range("E41").Value = "  +1.04%  "

# Row 42: PaxDollar
$ws.Range("E42").Value = "  -0.10%  "

# Row 43: Quant
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.95"
$ws.Range("E43").Value = "  +0.16%  "

# Row 44: RocketPoolETH
$ws.Range("D44").Value = "1.843.96"
$ws.Range("E44").Value = "  +0.75%  "

# Row 45: BabyDogeCoin
$ws.Range("E45").Value = "  +0.78%  "

# Row 46: Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.71"
$ws.Range("E46").Value = "  +1.19%  "

# Row 47: EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.198"
$ws.Range("E47").Value = "  +0.14%  "

# Row 48: Frax
$ws.Range("E48").Value = "  -0.37%  "

# Row 49: Cronos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05286"
$ws.Range("E49").Value = "  +0.24%  "

# Row 50: Mantle
$ws.Range("E50").Value = "  +0.19%  "

# Row 51: Aptos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.069"
$ws.Range("E51").Value = "  -0.34%  "
